$d = $word.ActiveDocument

# Locate the last paragraph in the document body (currently the lone empty
# paragraph right after the table, just before the final section break).
$n = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($n).Range

# Make room: four fresh empty paragraphs after the anchor. We will fill the
# middle two with the new reference-list text and drop the unused trailing
# spare at the end (Word always leaves one extra paragraph mark dangling
# when content is appended at the very end of the document body).
$null = $anchor.InsertParagraphAfter()
$null = $anchor.InsertParagraphAfter()
$null = $anchor.InsertParagraphAfter()
$null = $anchor.InsertParagraphAfter()

$n = $d.Paragraphs.Count

$emptyXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $d.Paragraphs.Item($n - 3).Range.InsertXML($emptyXml)

$namaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Nama </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>belakang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nama</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>depan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tahun</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Judul</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sumber</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Nama </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>perguruan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tinggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dan </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fakultas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $d.Paragraphs.Item($n - 2).Range.InsertXML($namaXml)

$citeXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="567" w:hanging="567"/><w:jc w:val="both"/></w:pPr><w:bookmarkStart w:id="1" w:name="_Hlk108944177"/><w:proofErr w:type="spellStart"/><w:r><w:t>Alfiansyah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rizky</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. 2021. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Implementasi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Single Exponential Smoothing Method</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sebagai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Dasar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pengendalian</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Persediaan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bahan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Baku Di </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Restoran</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cepat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Saji</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Fun Chicken </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tumpang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tugas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Akhir </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Skripsi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Program </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Studi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> S1 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sistem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Informasi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Universitas Merdeka Malang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Fakultas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Teknologi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Informasi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:bookmarkEnd w:id="1"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $d.Paragraphs.Item($n - 1).Range.InsertXML($citeXml)

# Remove the leftover spare paragraph mark left dangling at the very end.
$null = $d.Paragraphs.Item($d.Paragraphs.Count).Range.Delete()

Write-Output "Reference entry inserted. Paragraphs.Count is now $($d.Paragraphs.Count)."
